# chl-a spatial res in table 1
#
# Table 1 (Sheet1, A1:G9) lists environmental predictors. Two cells change:
#   E2 (Temporal Averaging for Sea Surface Temperature) : "2 days"   -> "10 days"
#   F5 (Spatial Averaging for Chlorophyll-a / chl-a)     : "Variable*" -> "4 km*"
#
# The author also left the selection on a different cell (F25) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "10 days"
$ws.Range("F5").Value = "4 km*"

# Reflect the saved selection state from the commit (best effort - the
# scroll position itself is cosmetic view state).
$ws.Range("F25").Select()
